$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.492932558059692
$ws.Range("B1").Value = 1.448063015937805
$ws.Range("C1").Value = 7.285219669342041
$ws.Range("D1").Value = 1.717133164405823
$ws.Range("E1").Value = 1.000354170799255
